# Apply cryptos list update (Mon Apr 10 09:42:20 UTC 2023)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.514.62"
$ws.Range("E2").Value = "'  -0.26%  "
$ws.Range("D3").Value = "'1.871.35"
$ws.Range("E3").Value = "'  -0.59%  "
$ws.Range("D4").Value = "'1.007"
$ws.Range("E4").Value = "'  -1.66%  "
$ws.Range("D5").Value = "'315.31"
$ws.Range("E5").Value = "'  -1.03%  "
$ws.Range("E6").Value = "'  -1.87%  "
$ws.Range("D7").Value = "'0.5093"
$ws.Range("E7").Value = "'  -1.15%  "
$ws.Range("D8").Value = "'0.3903"
$ws.Range("E8").Value = "'  -1.14%  "
$ws.Range("D9").Value = "'0.08357"
$ws.Range("E9").Value = "'  +0.04%  "
$ws.Range("D10").Value = "'1.106"
$ws.Range("E10").Value = "'  -1.30%  "
$ws.Range("D11").Value = "'6.215"
$ws.Range("E11").Value = "'  -0.84%  "
$ws.Range("D12").Value = "'1.873.66"
$ws.Range("E12").Value = "'  +0.87%  "
$ws.Range("D13").Value = "'20.38"
$ws.Range("E13").Value = "'  -0.70%  "
$ws.Range("D14").Value = "'7.282"
$ws.Range("E14").Value = "'  +0.22%  "
$ws.Range("E15").Value = "'  -1.64%  "
$ws.Range("E16").Value = "'  -1.05%  "
$ws.Range("D17").Value = "'91.11"
$ws.Range("E17").Value = "'  -0.40%  "
$ws.Range("D18").Value = "'0.06726"
$ws.Range("E18").Value = "'  -0.75%  "
$ws.Range("D19").Value = "'17.76"
$ws.Range("E19").Value = "'  +0.07%  "
$ws.Range("D20").Value = "'1.006"
$ws.Range("E20").Value = "'  -2.00%  "
$ws.Range("D21").Value = "'5.923"
$ws.Range("E21").Value = "'  -1.00%  "
$ws.Range("D22").Value = "'28.531.89"
$ws.Range("E22").Value = "'  -0.32%  "
$ws.Range("D23").Value = "'11.12"
$ws.Range("E23").Value = "'  -0.53%  "
$ws.Range("D24").Value = "'2.208"
$ws.Range("E24").Value = "'  -2.89%  "
$ws.Range("D25").Value = "'2.083.06"
$ws.Range("E25").Value = "'  +0.52%  "
$ws.Range("D26").Value = "'160.37"
$ws.Range("E26").Value = "'  -1.04%  "
$ws.Range("E27").Value = "'  -0.64%  "
$ws.Range("D28").Value = "'2.419"
$ws.Range("E28").Value = "'  +1.61%  "
$ws.Range("D29").Value = "'126.36"
$ws.Range("E29").Value = "'  -0.86%  "
$ws.Range("E30").Value = "'  -0.94%  "
$ws.Range("D31").Value = "'1.041"
$ws.Range("E31").Value = "'  +0.17%  "
$ws.Range("D32").Value = "'5.740"
$ws.Range("E32").Value = "'  -1.67%  "
$ws.Range("D33").Value = "'3.616"
$ws.Range("E33").Value = "'  -1.54%  "
$ws.Range("D34").Value = "'0.02457"
$ws.Range("E34").Value = "'  +0.21%  "
$ws.Range("D35").Value = "'0.06573"
$ws.Range("E35").Value = "'  +0.59%  "
$ws.Range("B36").Value = "'FraxShare"
$ws.Range("C36").Value = "'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D36").Value = "'8.915"
$ws.Range("E36").Value = "'  -2.95%  "
$ws.Range("B37").Value = "'Algorand"
$ws.Range("C37").Value = "'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D37").Value = "'0.2164"
$ws.Range("E37").Value = "'  -1.33%  "
$ws.Range("D38").Value = "'5.029"
$ws.Range("E38").Value = "'  +0.37%  "
$ws.Range("D39").Value = "'1.183"
$ws.Range("E39").Value = "'  -0.70%  "
$ws.Range("D40").Value = "'1.239"
$ws.Range("E40").Value = "'  -1.17%  "
$ws.Range("D41").Value = "'0.6373"
$ws.Range("E41").Value = "'  -1.47%  "
$ws.Range("D42").Value = "'11.10"
$ws.Range("E42").Value = "'  -1.09%  "
$ws.Range("D43").Value = "'1.006"
$ws.Range("E43").Value = "'  -1.84%  "
$ws.Range("D44").Value = "'0.6009"
$ws.Range("E44").Value = "'  -0.77%  "
$ws.Range("E45").Value = "'  +0.31%  "
$ws.Range("D46").Value = "'3.683"
$ws.Range("E46").Value = "'  -1.12%  "
$ws.Range("D47").Value = "'2.004"
$ws.Range("E47").Value = "'  +0.11%  "
$ws.Range("E48").Value = "'  -0.11%  "
$ws.Range("D49").Value = "'122.08"
$ws.Range("E49").Value = "'  -0.22%  "
$ws.Range("D50").Value = "'0.06812"
$ws.Range("E50").Value = "'  -0.97%  "
$ws.Range("D51").Value = "'76.40"
$ws.Range("E51").Value = "'  +0.07%  "
